$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newRows = @(
    @{ Row = 56; Date = "4/23/2021"; Url = "https://raw.githubusercontent.com/simonw/cdc-vaccination-history/6be51bd4348df57c1533fd1a13d3e0fcdd0107c7/states.json" },
    @{ Row = 57; Date = "4/24/2021"; Url = "https://raw.githubusercontent.com/simonw/cdc-vaccination-history/c9a530dc6d085617a9ea6d6669b4f9ef8ba3fd50/states.json" },
    @{ Row = 58; Date = "4/25/2021"; Url = "https://raw.githubusercontent.com/simonw/cdc-vaccination-history/518c8623cf6c257adbc938fe3ebcce965d2df854/states.json" }
)

foreach ($item in $newRows) {
    $r = $item.Row
    $ws.Cells.Item($r - 1, 1).Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)
    $ws.Cells.Item($r, 1).Value = $item.Date
    $ws.Cells.Item($r, 2).Value = $item.Url
}

$ws.Range("F49").Select()
